$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new bibliography entries (Microsoft IE11 system requirements, and
# Microsoft "What is User Account Control?") at the bottom of the existing
# reference table, then re-sort the whole table (A2:B20) by column B so the
# new rows land in their correct alphabetical position.
# ---------------------------------------------------------------------------

# New row 19: Microsoft. (n.d.-a). Internet Explorer system requirements IE11.
$ws.Range("A19").Value = "(Microsoft, n.d.-a)"

$b19 = $ws.Range("B19")
$b19.Value = "Microsoft. (n.d.-a). Internet Explorer system requirements IE11. Microsoft. Retrieved April 19, 2015, from http://windows.microsoft.com/en-au/internet-explorer/ie-system-requirements#ie=ie-11"
$b19.Characters(22,42).Font.Italic = $true
$b19.Characters(64,2).Font.Name = "Calibri"
$b19.Characters(66,9).Font.Italic = $true
$b19.Characters(75,117).Font.Name = "Calibri"
$b19.HorizontalAlignment = -4131
$b19.VerticalAlignment = -4108

# New row 20: Microsoft. (n.d.-b). What is User Account Control? - Windows Help.
$ws.Range("A20").Value = "(Microsoft, n.d.-b)"

$b20 = $ws.Range("B20")
$b20.Value = "Microsoft. (n.d.-b). What is User Account Control? - Windows Help. Retrieved March 14, 2015, from http://windows.microsoft.com/en-au/windows/what-is-user-account-control#1TC=windows-vista"
$b20.Characters(22,44).Font.Italic = $true
$b20.Characters(66,122).Font.Name = "Calibri"
$b20.HorizontalAlignment = -4131
$b20.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Re-sort the full reference table (excluding the header row) ascending by
# the reference column (B), exactly as "Data > Sort A to Z" would.
# ---------------------------------------------------------------------------
$sortRange = $ws.Range("A2:B20")
$sortKey = $ws.Range("B20")
$sortRange.Sort($sortKey, 1)

# Move the active selection, matching where the editor's cursor ended up.
$ws.Range("B9").Select()
